# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Periodo Mora" list (column E, rows 16-60) is re-sequenced from
# descending (2003 .. 1607) to ascending (1607 .. 2003); the "Valor Mora"
# (column F) follows the period it is attached to, and "Salario Basico"
# (column G) is refreshed to the new flat value for every period row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$periods = @(
    "1607","1608","1609","1610","1611","1612",
    "1701","1702","1703","1704","1705","1706","1707","1708","1709","1710","1711","1712",
    "1801","1802","1803","1804","1805","1806","1807","1808","1809","1810","1811","1812",
    "1901","1902","1903","1904","1905","1906","1907","1908","1909","1910","1911","1912",
    "2001","2002","2003"
)

$valorMora = @(
    27578,27578,27578,27578,27578,27578,
    27578,27578,27578,27578,27578,27578,27578,27578,27578,27578,27578,27578,
    27578,27578,27578,27578,27578,27578,27578,27578,31249,31249,31249,31249,
    31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,
    31249,31249,31249
)

$salarioBasico = 781242

$firstRow = 16
for ($i = 0; $i -lt $periods.Length; $i++) {
    $row = $firstRow + $i
    $ws.Cells.Item($row, 5).Value = $periods[$i]
    $ws.Cells.Item($row, 6).Value = $valorMora[$i]
    $ws.Cells.Item($row, 7).Value = $salarioBasico
}
